$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "asep"
$ws.Range("A5").Value = "siti"
$ws.Range("A6").Value = "aska"
$ws.Range("A7").Value = "saka"

$ws.Range("A8").Select()
